$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The run we need to split used to read:
#   "It depends on how many USB ports we have in our PC since we can
#    connect with one instrument with each USB port via the VISA
#    protocol."
# It must become four runs:
#   1) "It depends on how many USB ports we have in our PC"      (rFonts hint=eastAsia, unchanged rPr)
#   2) ", because "                                               (rFonts hint=default)
#   3) "we can connect with one instrument with each USB port
#       via the VISA protocol."                                  (rFonts hint=eastAsia)
#   4) " So, if there are X usb ports, we can connect X
#        instruments via the VISA protocol."                     (rFonts hint=default)
# ------------------------------------------------------------------

$oldText = "It depends on how many USB ports we have in our PC since we can connect with one instrument with each USB port via the VISA protocol."

$part1 = "It depends on how many USB ports we have in our PC"
$part2 = ", because "
$part3 = "we can connect with one instrument with each USB port via the VISA protocol."
$part4 = " So, if there are X usb ports, we can connect X instruments via the VISA protocol."

# Find the paragraph that still has the original (un-split) sentence.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq $oldText) {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate the paragraph containing the target sentence."
}

$paraStart = $targetPara.Range.Start

# Range covering exactly the run's text (excludes the paragraph mark).
$runRange = $d.Range($paraStart, $paraStart + $oldText.Length)

# The paragraph ends with a (hidden) "_GoBack" bookmark right after the
# run. Because our replacement range reaches all the way to the end of
# the run's text, that bookmark would otherwise get stretched so it
# wraps the freshly inserted runs instead of simply following them. To
# keep it in its original (collapsed, trailing) position we delete it
# first and then re-create it explicitly as part of the inserted XML,
# right after the last new run - matching the original layout.
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
    # no _GoBack bookmark present - nothing to do
}

function Esc-Xml($s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

$rPr1 = '<w:rPr><w:rFonts w:hint="eastAsia"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>'
$rPr2 = '<w:rPr><w:rFonts w:hint="default"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>'
$rPr3 = '<w:rPr><w:rFonts w:hint="eastAsia"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>'
$rPr4 = '<w:rPr><w:rFonts w:hint="default"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>'

$run1 = '<w:r>' + $rPr1 + '<w:t>' + (Esc-Xml $part1) + '</w:t></w:r>'
$run2 = '<w:r>' + $rPr2 + '<w:t xml:space="preserve">' + (Esc-Xml $part2) + '</w:t></w:r>'
$run3 = '<w:r>' + $rPr3 + '<w:t>' + (Esc-Xml $part3) + '</w:t></w:r>'
$run4 = '<w:r>' + $rPr4 + '<w:t xml:space="preserve">' + (Esc-Xml $part4) + '</w:t></w:r>'

$bookmarkXml = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' + $run1 + $run2 + $run3 + $run4 + $bookmarkXml + '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$runRange.InsertXML($xmlFrag)
